$wb = $excel.ActiveWorkbook

# Rename the "Include #0" sheet to "Include from RxNorm"
$wsInclude = $wb.Worksheets.Item("Include #0")
$wsInclude.Name = "Include from RxNorm"

# Update the Metadata sheet values (revert 0.1.8 merge back to 0.1.6 state)
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B3").Value = "0.1.6"
$wsMeta.Range("B6").Value = "active"
$wsMeta.Range("B8").Value = "2023-05-05T10:50:04-05:00"
$wsMeta.Range("B10").Value = "No display for ContactDetail"
$wsMeta.Range("B11").Value = "No display for ContactDetail"

# Remove the "Jurisdiction" row entirely (row 12), shifting later rows up
$wsMeta.Rows.Item(12).Delete()
